$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 89540105
$ws.Range("B2").Value = 9009
$ws.Range("E2").Value = 101603
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "Furustumpbagge"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "Plegaderus saucius"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "Erichson, 1834"
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "2"
$ws.Range("K2").NumberFormat = "@"
$ws.Range("K2").Value = "imago/adult"
$ws.Range("Q2").Value = 531886.6182379224
$ws.Range("R2").Value = 6623214.633886824
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value = "2020-11-25"
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value = "2020-11-25"
$ws.Range("AC2").NumberFormat = "@"
$ws.Range("AC2").Value = "inventering för Lst"
$ws.Range("AI2").NumberFormat = "@"
$ws.Range("AI2").Value = "äldre tallskog"
$ws.Range("AO2").NumberFormat = "@"
$ws.Range("AO2").Value = "stående nydöd tall"
$ws.Range("AW2").NumberFormat = "@"
$ws.Range("AW2").Value = "Olof Hedgren"
$ws.Range("AX2").NumberFormat = "@"
$ws.Range("AX2").Value = "Olof Hedgren"
$ws.Range("A3").Value = 98450400
$ws.Range("Q3").Value = 532384.3961978648
$ws.Range("R3").Value = 6623456.331745383
$ws.Range("A4").Value = 98450398
$ws.Range("Q4").Value = 532412.1361300815
$ws.Range("R4").Value = 6623401.632868396
$ws.Range("A5").Value = 98450410
$ws.Range("B5").Value = 90653
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "LC"
$ws.Range("E5").Value = 4364
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "Dropptaggsvamp"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "Hydnellum ferrugineum"
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q5").Value = 531846.2572265175
$ws.Range("R5").Value = 6623210.760434012
$ws.Range("A6").Value = 98450399
$ws.Range("Q6").Value = 532422.5980015037
$ws.Range("R6").Value = 6623417.853099325
$ws.Range("A7").Value = 98450403
$ws.Range("Q7").Value = 532144.1002281209
$ws.Range("R7").Value = 6623335.298360098
$ws.Range("A8").Value = 98450402
$ws.Range("B8").Value = 90676
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "NT"
$ws.Range("E8").Value = 5966
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "Motaggsvamp"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "Sarcodon squamosus"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H8").Value = "(Schaeff.) Quél."
$ws.Range("Q8").Value = 532224.744860352
$ws.Range("R8").Value = 6623351.620086531
$ws.Range("S8").Value = 10
$ws.Range("Y8").NumberFormat = "@"
$ws.Range("Y8").Value = "2021-10-08"
$ws.Range("Z8").NumberFormat = "@"
$ws.Range("Z8").Value = "00:00"
$ws.Range("AA8").NumberFormat = "@"
$ws.Range("AA8").Value = "2021-10-08"
$ws.Range("AB8").NumberFormat = "@"
$ws.Range("AB8").Value = "00:00"
$ws.Range("AI8").NumberFormat = "@"
$ws.Range("AI8").Value = "Äldre tallskog"
$ws.Range("AW8").NumberFormat = "@"
$ws.Range("AW8").Value = "Jacob Rudhe"
$ws.Range("AX8").NumberFormat = "@"
$ws.Range("AX8").Value = "Jacob Rudhe, Mårten Berglind"
$ws.Range("A9").Value = 98450401
$ws.Range("B9").Value = 90676
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "NT"
$ws.Range("E9").Value = 5966
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "Motaggsvamp"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "Sarcodon squamosus"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H9").Value = "(Schaeff.) Quél."
$ws.Range("Q9").Value = 532300.0021253217
$ws.Range("R9").Value = 6623465.176121892
$ws.Range("S9").Value = 10
$ws.Range("Y9").NumberFormat = "@"
$ws.Range("Y9").Value = "2021-10-08"
$ws.Range("Z9").NumberFormat = "@"
$ws.Range("Z9").Value = "00:00"
$ws.Range("AA9").NumberFormat = "@"
$ws.Range("AA9").Value = "2021-10-08"
$ws.Range("AB9").NumberFormat = "@"
$ws.Range("AB9").Value = "00:00"
$ws.Range("AI9").NumberFormat = "@"
$ws.Range("AI9").Value = "Äldre tallskog"
$ws.Range("AW9").NumberFormat = "@"
$ws.Range("AW9").Value = "Jacob Rudhe"
$ws.Range("AX9").NumberFormat = "@"
$ws.Range("AX9").Value = "Jacob Rudhe, Mårten Berglind"
$ws.Range("A10").Value = 98450406
$ws.Range("B10").Value = 90653
$ws.Range("Q10").Value = 531960.9745575936
$ws.Range("R10").Value = 6623258.115823878
$ws.Range("S10").Value = 10
$ws.Range("Y10").NumberFormat = "@"
$ws.Range("Y10").Value = "2021-10-08"
$ws.Range("Z10").NumberFormat = "@"
$ws.Range("Z10").Value = "00:00"
$ws.Range("AA10").NumberFormat = "@"
$ws.Range("AA10").Value = "2021-10-08"
$ws.Range("AB10").NumberFormat = "@"
$ws.Range("AB10").Value = "00:00"
$ws.Range("AI10").NumberFormat = "@"
$ws.Range("AI10").Value = "Äldre tallskog"
$ws.Range("AW10").NumberFormat = "@"
$ws.Range("AW10").Value = "Jacob Rudhe"
$ws.Range("AX10").NumberFormat = "@"
$ws.Range("AX10").Value = "Jacob Rudhe, Mårten Berglind"
$ws.Range("A11").Value = 98450405
$ws.Range("B11").Value = 90653
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "LC"
$ws.Range("E11").Value = 4364
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "Dropptaggsvamp"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "Hydnellum ferrugineum"
$ws.Range("H11").NumberFormat = "@"
$ws.Range("H11").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("I11").ClearContents()
$ws.Range("K11").ClearContents()
$ws.Range("Q11").Value = 532001.7303159089
$ws.Range("R11").Value = 6623333.569402254
$ws.Range("Y11").NumberFormat = "@"
$ws.Range("Y11").Value = "2021-10-08"
$ws.Range("AA11").NumberFormat = "@"
$ws.Range("AA11").Value = "2021-10-08"
$ws.Range("AC11").ClearContents()
$ws.Range("AI11").NumberFormat = "@"
$ws.Range("AI11").Value = "Äldre tallskog"
$ws.Range("AO11").ClearContents()
$ws.Range("AW11").NumberFormat = "@"
$ws.Range("AW11").Value = "Jacob Rudhe"
$ws.Range("AX11").NumberFormat = "@"
$ws.Range("AX11").Value = "Jacob Rudhe, Mårten Berglind"
$ws.Range("A12").Value = 98450404
$ws.Range("B12").Value = 90676
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 5966
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "Motaggsvamp"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "Sarcodon squamosus"
$ws.Range("H12").NumberFormat = "@"
$ws.Range("H12").Value = "(Schaeff.) Quél."
$ws.Range("Q12").Value = 532106.7513082939
$ws.Range("R12").Value = 6623274.994264505
$ws.Range("A13").Value = 98450407
$ws.Range("B13").Value = 90676
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "NT"
$ws.Range("E13").Value = 5966
$ws.Range("F13").NumberFormat = "@"
$ws.Range("F13").Value = "Motaggsvamp"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "Sarcodon squamosus"
$ws.Range("H13").NumberFormat = "@"
$ws.Range("H13").Value = "(Schaeff.) Quél."
$ws.Range("Q13").Value = 531962.2047132516
$ws.Range("R13").Value = 6623232.419310441
$ws.Range("A14").Value = 112379079
$ws.Range("B14").Value = 90814
$ws.Range("Q14").Value = 532226
$ws.Range("R14").Value = 6623334
$ws.Range("S14").Value = 25
$ws.Range("Y14").NumberFormat = "@"
$ws.Range("Y14").Value = "2023-08-27"
$ws.Range("Z14").ClearContents()
$ws.Range("AA14").NumberFormat = "@"
$ws.Range("AA14").Value = "2023-09-18"
$ws.Range("AB14").ClearContents()
$ws.Range("AI14").ClearContents()
$ws.Range("AW14").NumberFormat = "@"
$ws.Range("AW14").Value = "Mikael Hagström"
$ws.Range("AX14").NumberFormat = "@"
$ws.Range("AX14").Value = "Mikael Hagström"
$ws.Range("A15").Value = 112379171
$ws.Range("B15").Value = 90814
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "LC"
$ws.Range("E15").Value = 4364
$ws.Range("F15").NumberFormat = "@"
$ws.Range("F15").Value = "Dropptaggsvamp"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "Hydnellum ferrugineum"
$ws.Range("H15").NumberFormat = "@"
$ws.Range("H15").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q15").Value = 532169
$ws.Range("R15").Value = 6623457
$ws.Range("S15").Value = 25
$ws.Range("Y15").NumberFormat = "@"
$ws.Range("Y15").Value = "2023-08-27"
$ws.Range("Z15").ClearContents()
$ws.Range("AA15").NumberFormat = "@"
$ws.Range("AA15").Value = "2023-09-18"
$ws.Range("AB15").ClearContents()
$ws.Range("AI15").ClearContents()
$ws.Range("AW15").NumberFormat = "@"
$ws.Range("AW15").Value = "Mikael Hagström"
$ws.Range("AX15").NumberFormat = "@"
$ws.Range("AX15").Value = "Mikael Hagström"
$ws.Range("A16").Value = 112379172
$ws.Range("B16").Value = 90814
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "LC"
$ws.Range("E16").Value = 4364
$ws.Range("F16").NumberFormat = "@"
$ws.Range("F16").Value = "Dropptaggsvamp"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "Hydnellum ferrugineum"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q16").Value = 532294
$ws.Range("R16").Value = 6623516
$ws.Range("S16").Value = 25
$ws.Range("Y16").NumberFormat = "@"
$ws.Range("Y16").Value = "2023-08-27"
$ws.Range("Z16").ClearContents()
$ws.Range("AA16").NumberFormat = "@"
$ws.Range("AA16").Value = "2023-09-18"
$ws.Range("AB16").ClearContents()
$ws.Range("AI16").ClearContents()
$ws.Range("AW16").NumberFormat = "@"
$ws.Range("AW16").Value = "Mikael Hagström"
$ws.Range("AX16").NumberFormat = "@"
$ws.Range("AX16").Value = "Mikael Hagström"
